$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift error-series values left by one column (drop oldest quarter,
# add newly evaluated ifoCAST quarter where available)
$ws.Range("B2").Value = 0.2879090979994584
$ws.Range("C2").Value = -1.425880358899853
$ws.Range("D2").Value = 1.402475014699119
$ws.Range("E2").Value = -0.8328575851670005
$ws.Range("F2").Value = 0.5184774727506619
$ws.Range("G2").Value = -0.1029604570662399
$ws.Range("H2").Value = 0.3976002401245912
$ws.Range("I2").Value = -0.2703078322215502
$ws.Range("J2").Value = 0.1586931430164528
$ws.Range("K2").Value = 0.2163646915946629
$ws.Range("B3").Value = -1.443434480259818
$ws.Range("C3").Value = 1.384920893339154
$ws.Range("D3").Value = -0.8504117065269649
$ws.Range("E3").Value = 0.5009233513906975
$ws.Range("F3").Value = -0.1205145784262043
$ws.Range("G3").Value = 0.3800461187646267
$ws.Range("H3").Value = -0.2878619535815147
$ws.Range("I3").Value = 0.1411390216564884
$ws.Range("J3").Value = 0.1988105702346985
$ws.Range("K3").Value = 0.322776941072984
$ws.Range("B4").Value = 1.473028212290161
$ws.Range("C4").Value = -0.7623043875759586
$ws.Range("D4").Value = 0.5890306703417038
$ws.Range("E4").Value = -0.0324072594751981
$ws.Range("F4").Value = 0.4681534377156329
$ws.Range("G4").Value = -0.1997546346305085
$ws.Range("H4").Value = 0.2292463406074946
$ws.Range("I4").Value = 0.2869178891857047
$ws.Range("J4").Value = 0.4108842600239903
$ws.Range("K4").Value = -0.4825338632108016
$ws.Range("B5").Value = -0.5837297540881751
$ws.Range("C5").Value = 0.7676053038294873
$ws.Range("D5").Value = 0.1461673740125855
$ws.Range("E5").Value = 0.6467280712034165
$ws.Range("F5").Value = -0.02118000114272489
$ws.Range("G5").Value = 0.4078209740952782
$ws.Range("H5").Value = 0.4654925226734883
$ws.Range("I5").Value = 0.5894588935117738
$ws.Range("J5").Value = -0.303959229723018
$ws.Range("K5").Value = 0.4661714972207444
$ws.Range("B6").Value = 1.679632531582137
$ws.Range("C6").Value = 1.058194601765235
$ws.Range("D6").Value = 1.558755298956066
$ws.Range("E6").Value = 0.8908472266099251
$ws.Range("F6").Value = 1.319848201847928
$ws.Range("G6").Value = 1.377519750426138
$ws.Range("H6").Value = 1.501486121264424
$ws.Range("I6").Value = 0.608067998029632
$ws.Range("J6").Value = 1.378198724973394
$ws.Range("K6").Value = 1.11229800409388
$ws.Range("B7").Value = 0.1142203657994787
$ws.Range("C7").Value = 0.6147810629903097
$ws.Range("D7").Value = -0.0531270093558317
$ws.Range("E7").Value = 0.3758739658821714
$ws.Range("F7").Value = 0.4335455144603815
$ws.Range("G7").Value = 0.557511885298667
$ws.Range("H7").Value = -0.3359062379361248
$ws.Range("I7").Value = 0.4342244890076376
$ws.Range("J7").Value = 0.1683237681281231
$ws.Range("B8").Value = 0.6187489605034189
$ws.Range("C8").Value = -0.04915911184272259
$ws.Range("D8").Value = 0.3798418633952805
$ws.Range("E8").Value = 0.4375134119734906
$ws.Range("F8").Value = 0.5614797828117761
$ws.Range("G8").Value = -0.3319383404230157
$ws.Range("H8").Value = 0.4381923865207467
$ws.Range("I8").Value = 0.1722916656412322
$ws.Range("B9").Value = 0.08648097832751878
$ws.Range("C9").Value = 0.5154819535655218
$ws.Range("D9").Value = 0.573153502143732
$ws.Range("E9").Value = 0.6971198729820175
$ws.Range("F9").Value = -0.1962982502527744
$ws.Range("G9").Value = 0.5738324766909881
$ws.Range("H9").Value = 0.3079317558114735
$ws.Range("B10").Value = 0.2746757717098572
$ws.Range("C10").Value = 0.3323473202880673
$ws.Range("D10").Value = 0.4563136911263528
$ws.Range("E10").Value = -0.4371044321084391
$ws.Range("F10").Value = 0.3330262948353234
$ws.Range("G10").Value = 0.06712557395580883
$ws.Range("B11").Value = 0.2870161050359709
$ws.Range("C11").Value = 0.4109824758742565
$ws.Range("D11").Value = -0.4824356473605354
$ws.Range("E11").Value = 0.287695079583227
$ws.Range("F11").Value = 0.02179435870371246
$ws.Range("B12").Value = 0.3441210539382026
$ws.Range("C12").Value = -0.5492970692965893
$ws.Range("D12").Value = 0.2208336576471732
$ws.Range("E12").Value = -0.04506706323234141
$ws.Range("B13").Value = -0.5788832716533059
$ws.Range("C13").Value = 0.1912474552904566
$ws.Range("D13").Value = -0.07465326558905801
$ws.Range("B14").Value = 0.1730967985608157
$ws.Range("C14").Value = -0.0928039223186989
$ws.Range("B15").Value = -0.1108357465673982

# Clear trailing cells that no longer have data
$ws.Range("K7").ClearContents()
$ws.Range("J8").ClearContents()
$ws.Range("I9").ClearContents()
$ws.Range("H10").ClearContents()
$ws.Range("G11").ClearContents()
$ws.Range("F12").ClearContents()
$ws.Range("E13").ClearContents()
$ws.Range("D14").ClearContents()
$ws.Range("C15").ClearContents()
$ws.Range("B16").ClearContents()
